# Insert a new weekly data row at row 134 of the "Perejil" sheet (Vega Modelo
# de Temuco), pushing all existing rows 134-261 down to 135-262.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134:261 down one row, creating a blank row 134.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new data point.
$ws.Range("A134").Value = 10
$ws.Range("B134").Value = "Vega Modelo de Temuco"
$ws.Range("C134").Value = "La Araucanía"
$ws.Range("D134").Value = 44587
$ws.Range("E134").Value = 9
$ws.Range("F134").Value = 100112044
$ws.Range("G134").Value = "Perejil"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 55
$ws.Range("K134").Value = 5000
$ws.Range("L134").Value = 5000
$ws.Range("M134").Value = 5000
$ws.Range("N134").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O134").Value = "Provincia de Cautín"
$ws.Range("P134").Value = 1667
$ws.Range("Q134").Value = 3
$ws.Range("R134").Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D134").NumberFormat = $ws.Range("D135").NumberFormat()
